# Weekly data refresh: insert a new price observation at the top of the
# "Zapallo italiano" series (row 338), pushing the existing historical rows
# (338-360) down by one (339-361).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 338. This shifts rows
# 338..360 down to 339..361 (and all their formatting/values move with
# them), matching the rest of the sheet's existing row layout.
$ws.Rows.Item(338).Insert()

# Populate the newly inserted row 338 with the new weekly record.
$ws.Cells.Item(338, 1).Value  = 10
$ws.Cells.Item(338, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(338, 3).Value  = "La Araucanía"
$ws.Cells.Item(338, 4).Value  = 44610
$ws.Cells.Item(338, 5).Value  = 9
$ws.Cells.Item(338, 6).Value  = 100112032
$ws.Cells.Item(338, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(338, 8).Value  = "Sin especificar"
$ws.Cells.Item(338, 9).Value  = "Primera"
$ws.Cells.Item(338, 10).Value = 50
$ws.Cells.Item(338, 11).Value = 12000
$ws.Cells.Item(338, 12).Value = 12000
$ws.Cells.Item(338, 13).Value = 12000
$ws.Cells.Item(338, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(338, 15).Value = "Región del Maule"
$ws.Cells.Item(338, 16).Value = 200
$ws.Cells.Item(338, 17).Value = 60
$ws.Cells.Item(338, 18).Value = "Hortaliza"
